$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.855.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "'2.444.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.99%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'523.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'130.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'0.562"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("D9").Value = "'2.446.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.91%  "
$ws.Range("D10").Value = "'0.0973"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("D12").Value = "'4.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.76%  "
$ws.Range("D13").Value = "'0.322"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.38%  "
$ws.Range("D14").Value = "'2.874.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.95%  "
$ws.Range("D15").Value = "'57.740.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "'21.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("E17").Value = "  -1.99%  "
$ws.Range("D18").Value = "'2.440.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.00%  "
$ws.Range("D19").Value = "'10.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.95%  "
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").Value = "'314.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("D22").Value = "'6.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'64.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").Value = "'0.405"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("D28").Value = "'7.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.73%  "
$ws.Range("D29").Value = "'173.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("D30").Value = "'0.0₃0732"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("D31").Value = "'1.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").Value = "'6.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("D33").Value = "'1.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.97%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").Value = "'17.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").Value = "'1.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.61%  "
$ws.Range("D38").Value = "'3.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.72%  "
$ws.Range("D39").Value = "'36.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").Value = "'1.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").Value = "'0.790"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("D42").Value = "'3.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.84%  "
$ws.Range("D43").Value = "'263.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.581"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.92%  "
$ws.Range("D45").Value = "'4.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.00%  "
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").Value = "'121.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("D48").Value = "'0.0491"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.14%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0210"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'16.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.51%  "
$ws.Range("D51").Value = "'16.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.62%  "
